$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M2").Value = 5.273684
$ws.Range("N2").Value = 15.821052
$ws.Range("O2").Value = 0.0510821201937383
$ws.Range("P2").Value = 0.0510821201937383
$ws.Range("Q2").Value = 1.086568756624
$ws.Range("R2").Value = 9.779118809616
$ws.Range("S2").Value = 0.02367655104596258
$ws.Range("T2").Value = 0.02367655104596259
$ws.Range("O3").Value = 0.5598845502029881
$ws.Range("P3").Value = 0.5598845502029881
$ws.Range("S3").Value = 0.259506361177855
$ws.Range("T3").Value = 0.259506361177855
$ws.Range("M4").Value = 32.95839133333334
$ws.Range("N4").Value = 98.87517400000002
$ws.Range("O4").Value = 0.3192425840231603
$ws.Range("P4").Value = 0.3192425840231604
$ws.Range("Q4").Value = 6.790615116754668
$ws.Range("R4").Value = 61.11553605079201
$ws.Range("S4").Value = 0.1479688648004844
$ws.Range("T4").Value = 0.1479688648004844
$ws.Range("M5").Value = 7.205150000000001
$ws.Range("N5").Value = 21.61545
$ws.Range("O5").Value = 0.06979074558011317
$ws.Range("P5").Value = 0.06979074558011318
$ws.Range("Q5").Value = 1.4845202854
$ws.Range("R5").Value = 13.3606825686
$ws.Range("S5").Value = 0.03234799464071365
$ws.Range("T5").Value = 0.03234799464071366
$ws.Range("G6").Value = 0.2384863333333333
$ws.Range("H6").Value = 0.715459
$ws.Range("I6").Value = 0.5365002283349842
$ws.Range("J6").Value = 0.5365002283349842
$ws.Range("M6").Value = 5.273684
$ws.Range("N6").Value = 15.821052
$ws.Range("O6").Value = 0.0510821201937383
$ws.Range("P6").Value = 0.0510821201937383
$ws.Range("Q6").Value = 1.257701560318667
$ws.Range("R6").Value = 11.319314042868
$ws.Range("S6").Value = 0.02740556914777571
$ws.Range("T6").Value = 0.02740556914777571
$ws.Range("G7").Value = 0.2384863333333333
$ws.Range("H7").Value = 0.715459
$ws.Range("I7").Value = 0.5365002283349842
$ws.Range("J7").Value = 0.5365002283349842
$ws.Range("O7").Value = 0.5598845502029881
$ws.Range("P7").Value = 0.5598845502029881
$ws.Range("Q7").Value = 13.78501263686644
$ws.Range("S7").Value = 0.3003781890251331
$ws.Range("T7").Value = 0.3003781890251331
$ws.Range("G8").Value = 0.2384863333333333
$ws.Range("H8").Value = 0.715459
$ws.Range("I8").Value = 0.5365002283349842
$ws.Range("J8").Value = 0.5365002283349842
$ws.Range("M8").Value = 32.95839133333334
$ws.Range("N8").Value = 98.87517400000002
$ws.Range("O8").Value = 0.3192425840231603
$ws.Range("P8").Value = 0.3192425840231604
$ws.Range("Q8").Value = 7.860125901651779
$ws.Range("R8").Value = 70.741133114866
$ws.Range("S8").Value = 0.1712737192226759
$ws.Range("T8").Value = 0.1712737192226759
$ws.Range("G9").Value = 0.2384863333333333
$ws.Range("H9").Value = 0.715459
$ws.Range("I9").Value = 0.5365002283349842
$ws.Range("J9").Value = 0.5365002283349842
$ws.Range("M9").Value = 7.205150000000001
$ws.Range("N9").Value = 21.61545
$ws.Range("O9").Value = 0.06979074558011317
$ws.Range("P9").Value = 0.06979074558011318
$ws.Range("Q9").Value = 1.718329804616667
$ws.Range("R9").Value = 15.46496824155
$ws.Range("S9").Value = 0.03744275093939951
$ws.Range("T9").Value = 0.03744275093939951
